# Generate Report for Handback
# Refresh the handback status report timestamps / priority for the new run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the first file
# moves forward from 12:15:27 to 12:16:19.
$wsOverview.Range("G2").Value = "2016-08-26 12:16:19"
$wsOverview.Range("G3").Value = "2016-08-26 12:16:19"

# Priority goes from "ht" (human translation) to "mt" (machine translation)
# for the first two files, on both the zh-cn and de-de sheets.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# zh-cn sheet: handoff/handback timestamps advance.
$wsZhCn.Range("H2").Value = "2016-08-26 12:16:15"
$wsZhCn.Range("H3").Value = "2016-08-26 12:16:15"
$wsZhCn.Range("K2").Value = "2016-08-26 12:16:32"
$wsZhCn.Range("K3").Value = "2016-08-26 12:16:32"

# de-de sheet: handoff datetime matches the Overview's regenerated timestamp,
# and the handback datetime also moves forward.
$wsDeDe.Range("H2").Value = "2016-08-26 12:16:19"
$wsDeDe.Range("H3").Value = "2016-08-26 12:16:19"
$wsDeDe.Range("K2").Value = "2016-08-26 12:16:38"
$wsDeDe.Range("K3").Value = "2016-08-26 12:16:38"
